# "TODAS LAS INTERFACES CON NUEVO DISENIO" -- update user access levels
# ("Nivel") on the Usuarios sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (JulianGuardian / BloodSlayer): Nivel 1 -> 3
$ws.Range("E3").Value = 3

# Row 6 (David Gonzalo Cordon Fontecha / The_Cortux): Nivel 3 -> 2
$ws.Range("E6").Value = 2

# Leave the cursor where the editor last clicked, one row below the table.
$ws.Range("D7").Select() | Out-Null
